$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1348.875
$ws.Range("I32").Value = 1380
$ws.Range("J32").Value = 1334.7273
$ws.Range("K32").Value = 1380
$ws.Range("L32").Value = 1334.7273
$ws.Range("M32").Value = -1054
$ws.Range("N32").Value = -1986.7273
$ws.Range("H64").Value = 5195.647
$ws.Range("I64").Value = 3193.9167
$ws.Range("J64").Value = 9999.799999999999
$ws.Range("K64").Value = 3193.9167
$ws.Range("L64").Value = 9999.799999999999
$ws.Range("M64").Value = -2945.9167
$ws.Range("N64").Value = -10495.8
$ws.Range("H67").Value = 5195.647
$ws.Range("I67").Value = 3193.9167
$ws.Range("J67").Value = 9999.799999999999
$ws.Range("K67").Value = 3193.9167
$ws.Range("L67").Value = 9999.799999999999
$ws.Range("M67").Value = -2335.9167
$ws.Range("N67").Value = -11715.8
$ws.Range("H82").Value = 10000
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 10000
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 30000
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -30812
$ws.Range("H85").Value = 10000
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 10000
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 30000
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -32808
$ws.Range("H98").Value = 1122.375
$ws.Range("I98").Value = 1122.375
$ws.Range("K98").Value = 1122.375
$ws.Range("M98").Value = 375.625
$ws.Range("H113").Value = 13609.3
$ws.Range("I113").Value = 16267.143
$ws.Range("K113").Value = 16267.143
$ws.Range("M113").Value = -13013.143
$ws.Range("H122").Value = 1122.375
$ws.Range("I122").Value = 1122.375
$ws.Range("K122").Value = 3367.125
$ws.Range("M122").Value = -917.125
$ws.Range("H132").Value = 1744.4762
$ws.Range("I132").Value = 1455.25
$ws.Range("K132").Value = 4365.75
$ws.Range("M132").Value = -1835.75
$ws.Range("H141").Value = 614.381
$ws.Range("I141").Value = 614.381
$ws.Range("K141").Value = 1843.143
$ws.Range("M141").Value = 3336.857

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3680.818
$ws.Range("I45").Value = 1635.6
$ws.Range("K45").Value = 1635.6
$ws.Range("M45").Value = -1258.6
$ws.Range("H61").Value = 900
$ws.Range("I61").Value = 900
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 900
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -688
$ws.Range("N61").ClearContents()
$ws.Range("H102").Value = 2218.6667
$ws.Range("I102").Value = 2218.6667
$ws.Range("K102").Value = 2218.6667
$ws.Range("M102").Value = -596.6667000000002
$ws.Range("H122").Value = 3346
$ws.Range("J122").Value = 4459
$ws.Range("L122").Value = 13377
$ws.Range("N122").Value = -18277
$ws.Range("H132").Value = 961
$ws.Range("I132").Value = 961
$ws.Range("K132").Value = 2883
$ws.Range("M132").Value = -353
$ws.Range("H135").Value = 30000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 30000
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 30000
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -40140
$ws.Range("H136").Value = 900
$ws.Range("I136").Value = 900
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2700
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -150
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2677.1765
$ws.Range("I105").Value = 2750.9285
$ws.Range("K105").Value = 2750.9285
$ws.Range("M105").Value = -1003.9285
$ws.Range("H134").Value = 1565.7826
$ws.Range("I134").Value = 1619.762
$ws.Range("J134").Value = 999
$ws.Range("K134").Value = 4859.286
$ws.Range("L134").Value = 2997
$ws.Range("M134").Value = -2324.286
$ws.Range("N134").Value = -8067
$ws.Range("H135").Value = 46666.332
$ws.Range("J135").Value = 49999.5
$ws.Range("L135").Value = 49999.5
$ws.Range("N135").Value = -60139.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2269.8333
$ws.Range("I31").Value = 1723.8
$ws.Range("K31").Value = 1723.8
$ws.Range("M31").Value = -1428.8
$ws.Range("H34").Value = 2269.8333
$ws.Range("I34").Value = 1723.8
$ws.Range("K34").Value = 1723.8
$ws.Range("M34").Value = -1521.8
$ws.Range("H107").Value = 549.125
$ws.Range("I107").Value = 549.125
$ws.Range("K107").Value = 549.125
$ws.Range("M107").Value = 1370.875
$ws.Range("H132").Value = 4587.6924
$ws.Range("I132").Value = 4587.6924
$ws.Range("K132").Value = 13763.0772
$ws.Range("M132").Value = -11233.0772
$ws.Range("H134").Value = 2268.5
$ws.Range("I134").Value = 1913.5555
$ws.Range("J134").Value = 3333.3333
$ws.Range("K134").Value = 5740.666499999999
$ws.Range("L134").Value = 9999.999899999999
$ws.Range("M134").Value = -3205.666499999999
$ws.Range("N134").Value = -15069.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1999.6666
$ws.Range("I68").Value = 1999
$ws.Range("K68").Value = 5997
$ws.Range("M68").Value = -5186
$ws.Range("H71").Value = 1999.6666
$ws.Range("I71").Value = 1999
$ws.Range("K71").Value = 17991
$ws.Range("M71").Value = -13935
$ws.Range("H86").Value = 1674.375
$ws.Range("I86").Value = 1674.375
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 5023.125
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -3837.125
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 1674.375
$ws.Range("I89").Value = 1674.375
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 15069.375
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -9141.375
$ws.Range("N89").ClearContents()
$ws.Range("H97").Value = 1525.4
$ws.Range("I97").Value = 1188.75
$ws.Range("J97").Value = 1749.8334
$ws.Range("K97").Value = 3566.25
$ws.Range("L97").Value = 5249.5002
$ws.Range("M97").Value = -3070.25
$ws.Range("N97").Value = -6241.5002
$ws.Range("H134").Value = 250000160
$ws.Range("I134").Value = 250000160
$ws.Range("K134").Value = 750000480
$ws.Range("M134").Value = -749995410

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H54").Value = 6678360.5
$ws.Range("I54").Value = 20000086
$ws.Range("J54").Value = 17497.5
$ws.Range("K54").Value = 20000086
$ws.Range("L54").Value = 17497.5
$ws.Range("M54").Value = -19999696
$ws.Range("N54").Value = -18277.5
$ws.Range("H122").Value = 3251.077
$ws.Range("I122").Value = 3397.818
$ws.Range("K122").Value = 10193.454
$ws.Range("M122").Value = -7743.454000000002
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H132").Value = 2618.2856
$ws.Range("I132").Value = 2243.1538
$ws.Range("K132").Value = 6729.4614
$ws.Range("M132").Value = -4199.4614

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 3330
$ws.Range("I13").Value = 2000
$ws.Range("J13").Value = 3995
$ws.Range("K13").Value = 2000
$ws.Range("L13").Value = 3995
$ws.Range("M13").Value = -1860
$ws.Range("N13").Value = -4275
$ws.Range("H16").Value = 1170.7142
$ws.Range("I16").Value = 1170.6
$ws.Range("J16").Value = 1171
$ws.Range("K16").Value = 1170.6
$ws.Range("L16").Value = 1171
$ws.Range("M16").Value = -1000.6
$ws.Range("N16").Value = -1511
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").ClearContents()
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()
$ws.Range("H53").Value = 15789
$ws.Range("I53").Value = 15789
$ws.Range("K53").Value = 15789
$ws.Range("M53").Value = -15271
$ws.Range("H55").Value = 564.2857
$ws.Range("I55").Value = 66.666664
$ws.Range("J55").Value = 937.5
$ws.Range("K55").Value = 66.666664
$ws.Range("L55").Value = 937.5
$ws.Range("M55").Value = 106.333336
$ws.Range("N55").Value = -1283.5
$ws.Range("H61").Value = 1770.75
$ws.Range("I61").Value = 1770.75
$ws.Range("K61").Value = 1770.75
$ws.Range("M61").Value = -1568.75
$ws.Range("H93").Value = 1822.75
$ws.Range("I93").Value = 1430.3334
$ws.Range("K93").Value = 1430.3334
$ws.Range("M93").Value = -182.3334
$ws.Range("H113").Value = 1770.75
$ws.Range("I113").Value = 1770.75
$ws.Range("K113").Value = 1770.75
$ws.Range("M113").Value = 399.25
$ws.Range("H132").Value = 1208
$ws.Range("I132").Value = 1119.6
$ws.Range("K132").Value = 3358.8
$ws.Range("M132").Value = -828.7999999999997

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("H80").Value = 29999
$ws.Range("J80").Value = 29999
$ws.Range("L80").Value = 29999
$ws.Range("N80").Value = -31995
$ws.Range("H83").Value = 29999
$ws.Range("J83").Value = 29999
$ws.Range("L83").Value = 89997
$ws.Range("N83").Value = -99981
$ws.Range("H100").Value = 6339684
$ws.Range("I100").Value = 8715316
$ws.Range("K100").Value = 17430632
$ws.Range("M100").Value = -17430091
$ws.Range("H122").Value = 3750
$ws.Range("I122").Value = 3750
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 11250
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -8800
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 2608.1428
$ws.Range("I126").Value = 2608.1428
$ws.Range("K126").Value = 7824.428400000001
$ws.Range("M126").Value = -5354.428400000001

